$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FG")

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43100
$ws.Range("F7").Value = 42735
$ws.Range("G7").Value = "NA"
$ws.Range("H7").Value = "NA"
$ws.Range("I7").Value = "NA"
$ws.Range("J7").Value = "NA"
$ws.Range("K7").Value = "NA"
$ws.Range("L7").Value = $null
$ws.Range("D8").Value = 735000
$ws.Range("E8").Value = 2079000
$ws.Range("F8").Value = 1525000
$ws.Range("G8").Value = "NA"
$ws.Range("H8").Value = "NA"
$ws.Range("I8").Value = "NA"
$ws.Range("J8").Value = "NA"
$ws.Range("K8").Value = "NA"
$ws.Range("L8").Value = $null
$ws.Range("D9").Value = 423000
$ws.Range("E9").Value = 1194000
$ws.Range("F9").Value = 811000
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("L9").Value = $null
$ws.Range("D10").Value = 312000
$ws.Range("E10").Value = 885000
$ws.Range("F10").Value = 714000
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "NA"
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("L10").Value = $null
$ws.Range("D11").Value = $null
$ws.Range("E11").Value = $null
$ws.Range("F11").Value = $null
$ws.Range("G11").Value = $null
$ws.Range("H11").Value = $null
$ws.Range("I11").Value = $null
$ws.Range("J11").Value = $null
$ws.Range("K11").Value = $null
$ws.Range("L11").Value = $null
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = $null
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = $null
$ws.Range("D14").Value = 55000
$ws.Range("E14").Value = 48000
$ws.Range("F14").Value = 61000
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
$ws.Range("K14").Value = "NA"
$ws.Range("L14").Value = $null
$ws.Range("D15").Value = 49000
$ws.Range("E15").Value = 233000
$ws.Range("F15").Value = 177000
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = "NA"
$ws.Range("I15").Value = "NA"
$ws.Range("J15").Value = "NA"
$ws.Range("K15").Value = "NA"
$ws.Range("L15").Value = $null
$ws.Range("D16").Value = $null
$ws.Range("E16").Value = $null
$ws.Range("F16").Value = $null
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = $null
$ws.Range("K16").Value = $null
$ws.Range("L16").Value = $null
$ws.Range("D17").Value = 677000
$ws.Range("E17").Value = 1653000
$ws.Range("F17").Value = 1181000
$ws.Range("G17").Value = "NA"
$ws.Range("H17").Value = "NA"
$ws.Range("I17").Value = "NA"
$ws.Range("J17").Value = "NA"
$ws.Range("K17").Value = "NA"
$ws.Range("L17").Value = $null
$ws.Range("D18").Value = 58000
$ws.Range("E18").Value = 426000
$ws.Range("F18").Value = 344000
$ws.Range("G18").Value = "NA"
$ws.Range("H18").Value = "NA"
$ws.Range("I18").Value = "NA"
$ws.Range("J18").Value = "NA"
$ws.Range("K18").Value = "NA"
$ws.Range("L18").Value = $null
$ws.Range("D19").Value = $null
$ws.Range("E19").Value = $null
$ws.Range("F19").Value = $null
$ws.Range("G19").Value = $null
$ws.Range("H19").Value = $null
$ws.Range("I19").Value = $null
$ws.Range("J19").Value = $null
$ws.Range("K19").Value = $null
$ws.Range("L19").Value = $null
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = "NA"
$ws.Range("H20").Value = "NA"
$ws.Range("I20").Value = "NA"
$ws.Range("J20").Value = "NA"
$ws.Range("K20").Value = "NA"
$ws.Range("L20").Value = $null
$ws.Range("D21").Value = 101000
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "NA"
$ws.Range("H21").Value = "NA"
$ws.Range("I21").Value = "NA"
$ws.Range("J21").Value = "NA"
$ws.Range("K21").Value = "NA"
$ws.Range("L21").Value = $null
$ws.Range("D22").Value = 29000
$ws.Range("E22").Value = 30000
$ws.Range("F22").Value = 28000
$ws.Range("G22").Value = "NA"
$ws.Range("H22").Value = "NA"
$ws.Range("I22").Value = "NA"
$ws.Range("J22").Value = "NA"
$ws.Range("K22").Value = "NA"
$ws.Range("L22").Value = $null
$ws.Range("D23").Value = 29000
$ws.Range("E23").Value = 396000
$ws.Range("F23").Value = 316000
$ws.Range("G23").Value = "NA"
$ws.Range("H23").Value = "NA"
$ws.Range("I23").Value = "NA"
$ws.Range("J23").Value = "NA"
$ws.Range("K23").Value = "NA"
$ws.Range("L23").Value = $null
$ws.Range("D24").Value = 16000
$ws.Range("E24").Value = 105000
$ws.Range("F24").Value = 111000
$ws.Range("G24").Value = "NA"
$ws.Range("H24").Value = "NA"
$ws.Range("I24").Value = "NA"
$ws.Range("J24").Value = "NA"
$ws.Range("K24").Value = "NA"
$ws.Range("L24").Value = $null
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = $null
$ws.Range("D26").Value = 13000
$ws.Range("E26").Value = 291000
$ws.Range("F26").Value = 205000
$ws.Range("G26").Value = "NA"
$ws.Range("H26").Value = "NA"
$ws.Range("I26").Value = "NA"
$ws.Range("J26").Value = "NA"
$ws.Range("K26").Value = "NA"
$ws.Range("L26").Value = $null
$ws.Range("D27").Value = -16000
$ws.Range("E27").Value = 289000
$ws.Range("F27").Value = 205000
$ws.Range("G27").Value = "NA"
$ws.Range("H27").Value = "NA"
$ws.Range("I27").Value = "NA"
$ws.Range("J27").Value = "NA"
$ws.Range("K27").Value = "NA"
$ws.Range("L27").Value = $null
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = $null
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = -131000
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = "NA"
$ws.Range("L29").Value = $null
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = $null
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = $null
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = "NA"
$ws.Range("H32").Value = "NA"
$ws.Range("I32").Value = "NA"
$ws.Range("J32").Value = "NA"
$ws.Range("K32").Value = "NA"
$ws.Range("L32").Value = $null
$ws.Range("D33").Value = -16000
$ws.Range("E33").Value = 158000
$ws.Range("F33").Value = 205000
$ws.Range("G33").Value = "NA"
$ws.Range("H33").Value = "NA"
$ws.Range("I33").Value = "NA"
$ws.Range("J33").Value = "NA"
$ws.Range("K33").Value = "NA"
$ws.Range("L33").Value = $null
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = $null
$ws.Range("D35").Value = -16000
$ws.Range("E35").Value = 158000
$ws.Range("F35").Value = 205000
$ws.Range("G35").Value = "NA"
$ws.Range("H35").Value = "NA"
$ws.Range("I35").Value = "NA"
$ws.Range("J35").Value = "NA"
$ws.Range("K35").Value = "NA"
$ws.Range("L35").Value = $null
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43100
$ws.Range("F38").Value = 42735
$ws.Range("G38").Value = "NA"
$ws.Range("H38").Value = "NA"
$ws.Range("I38").Value = "NA"
$ws.Range("J38").Value = "NA"
$ws.Range("K38").Value = "NA"
$ws.Range("L38").Value = $null
$ws.Range("D39").Value = $null
$ws.Range("E39").Value = $null
$ws.Range("F39").Value = $null
$ws.Range("G39").Value = $null
$ws.Range("H39").Value = $null
$ws.Range("I39").Value = $null
$ws.Range("J39").Value = $null
$ws.Range("K39").Value = $null
$ws.Range("L39").Value = $null
$ws.Range("D40").Value = $null
$ws.Range("E40").Value = $null
$ws.Range("F40").Value = $null
$ws.Range("G40").Value = $null
$ws.Range("H40").Value = $null
$ws.Range("I40").Value = $null
$ws.Range("J40").Value = $null
$ws.Range("K40").Value = $null
$ws.Range("L40").Value = $null
$ws.Range("D41").Value = 571000
$ws.Range("E41").Value = 1215000
$ws.Range("F41").Value = 1000
$ws.Range("G41").Value = "NA"
$ws.Range("H41").Value = "NA"
$ws.Range("I41").Value = "NA"
$ws.Range("J41").Value = "NA"
$ws.Range("K41").Value = "NA"
$ws.Range("L41").Value = $null
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = $null
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = $null
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = $null
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = 100
$ws.Range("G45").Value = "NA"
$ws.Range("H45").Value = "NA"
$ws.Range("I45").Value = "NA"
$ws.Range("J45").Value = "NA"
$ws.Range("K45").Value = "NA"
$ws.Range("L45").Value = $null
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = $null
$ws.Range("D47").Value = 23820000
$ws.Range("E47").Value = 23112000
$ws.Range("F47").Value = 690900
$ws.Range("G47").Value = "NA"
$ws.Range("H47").Value = "NA"
$ws.Range("I47").Value = "NA"
$ws.Range("J47").Value = "NA"
$ws.Range("K47").Value = "NA"
$ws.Range("L47").Value = $null
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = $null
$ws.Range("D49").Value = 1826000
$ws.Range("E49").Value = 1320000
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "NA"
$ws.Range("H49").Value = "NA"
$ws.Range("I49").Value = "NA"
$ws.Range("J49").Value = "NA"
$ws.Range("K49").Value = "NA"
$ws.Range("L49").Value = $null
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = $null
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = $null
$ws.Range("D52").Value = 343000
$ws.Range("E52").Value = 182000
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "NA"
$ws.Range("H52").Value = "NA"
$ws.Range("I52").Value = "NA"
$ws.Range("J52").Value = "NA"
$ws.Range("K52").Value = "NA"
$ws.Range("L52").Value = $null
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = $null
$ws.Range("D54").Value = 30945000
$ws.Range("E54").Value = 29923000
$ws.Range("F54").Value = 692100
$ws.Range("G54").Value = "NA"
$ws.Range("H54").Value = "NA"
$ws.Range("I54").Value = "NA"
$ws.Range("J54").Value = "NA"
$ws.Range("K54").Value = "NA"
$ws.Range("L54").Value = $null
$ws.Range("D55").Value = $null
$ws.Range("E55").Value = $null
$ws.Range("F55").Value = $null
$ws.Range("G55").Value = $null
$ws.Range("H55").Value = $null
$ws.Range("I55").Value = $null
$ws.Range("J55").Value = $null
$ws.Range("K55").Value = $null
$ws.Range("L55").Value = $null
$ws.Range("D56").Value = $null
$ws.Range("E56").Value = $null
$ws.Range("F56").Value = $null
$ws.Range("G56").Value = $null
$ws.Range("H56").Value = $null
$ws.Range("I56").Value = $null
$ws.Range("J56").Value = $null
$ws.Range("K56").Value = $null
$ws.Range("L56").Value = $null
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = 600
$ws.Range("G57").Value = "NA"
$ws.Range("H57").Value = "NA"
$ws.Range("I57").Value = "NA"
$ws.Range("J57").Value = "NA"
$ws.Range("K57").Value = "NA"
$ws.Range("L57").Value = $null
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = $null
$ws.Range("D59").Value = 28858000
$ws.Range("E59").Value = 26626000
$ws.Range("F59").Value = 300
$ws.Range("G59").Value = "NA"
$ws.Range("H59").Value = "NA"
$ws.Range("I59").Value = "NA"
$ws.Range("J59").Value = "NA"
$ws.Range("K59").Value = "NA"
$ws.Range("L59").Value = $null
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = $null
$ws.Range("D61").Value = 541000
$ws.Range("E61").Value = 412000
$ws.Range("F61").Value = 640700
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = $null
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = $null
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = $null
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = $null
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = $null
$ws.Range("D66").Value = 30055000
$ws.Range("E66").Value = 27960000
$ws.Range("F66").Value = 687100
$ws.Range("G66").Value = "NA"
$ws.Range("H66").Value = "NA"
$ws.Range("I66").Value = "NA"
$ws.Range("J66").Value = "NA"
$ws.Range("K66").Value = "NA"
$ws.Range("L66").Value = $null
$ws.Range("D67").Value = $null
$ws.Range("E67").Value = $null
$ws.Range("F67").Value = $null
$ws.Range("G67").Value = $null
$ws.Range("H67").Value = $null
$ws.Range("I67").Value = $null
$ws.Range("J67").Value = $null
$ws.Range("K67").Value = $null
$ws.Range("L67").Value = $null
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = $null
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = $null
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = $null
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = $null
$ws.Range("D72").Value = -167000
$ws.Range("E72").Value = -149000
$ws.Range("F72").Value = -300
$ws.Range("G72").Value = "NA"
$ws.Range("H72").Value = "NA"
$ws.Range("I72").Value = "NA"
$ws.Range("J72").Value = "NA"
$ws.Range("K72").Value = "NA"
$ws.Range("L72").Value = $null
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = $null
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = $null
$ws.Range("D76").Value = 890000
$ws.Range("E76").Value = 1963000
$ws.Range("F76").Value = 5000
$ws.Range("G76").Value = "NA"
$ws.Range("H76").Value = "NA"
$ws.Range("I76").Value = "NA"
$ws.Range("J76").Value = "NA"
$ws.Range("K76").Value = "NA"
$ws.Range("L76").Value = $null
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43100
$ws.Range("F80").Value = 42735
$ws.Range("G80").Value = "NA"
$ws.Range("H80").Value = "NA"
$ws.Range("I80").Value = "NA"
$ws.Range("J80").Value = "NA"
$ws.Range("K80").Value = "NA"
$ws.Range("L80").Value = $null
$ws.Range("D81").Value = -16000
$ws.Range("E81").Value = 158000
$ws.Range("F81").Value = 205000
$ws.Range("G81").Value = "NA"
$ws.Range("H81").Value = "NA"
$ws.Range("I81").Value = "NA"
$ws.Range("J81").Value = "NA"
$ws.Range("K81").Value = "NA"
$ws.Range("L81").Value = $null
$ws.Range("D82").Value = $null
$ws.Range("E82").Value = $null
$ws.Range("F82").Value = $null
$ws.Range("G82").Value = $null
$ws.Range("H82").Value = $null
$ws.Range("I82").Value = $null
$ws.Range("J82").Value = $null
$ws.Range("K82").Value = $null
$ws.Range("L82").Value = $null
$ws.Range("D83").Value = 43000
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "NA"
$ws.Range("H83").Value = "NA"
$ws.Range("I83").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = $null
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = $null
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = $null
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = $null
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = $null
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = $null
$ws.Range("D89").Value = 897000
$ws.Range("E89").Value = 401000
$ws.Range("F89").Value = 437000
$ws.Range("G89").Value = "NA"
$ws.Range("H89").Value = "NA"
$ws.Range("I89").Value = "NA"
$ws.Range("J89").Value = "NA"
$ws.Range("K89").Value = "NA"
$ws.Range("L89").Value = $null
$ws.Range("D90").Value = $null
$ws.Range("E90").Value = $null
$ws.Range("F90").Value = $null
$ws.Range("G90").Value = $null
$ws.Range("H90").Value = $null
$ws.Range("I90").Value = $null
$ws.Range("J90").Value = $null
$ws.Range("K90").Value = $null
$ws.Range("L90").Value = $null
$ws.Range("D91").Value = -7000
$ws.Range("E91").Value = -6000
$ws.Range("F91").Value = -10000
$ws.Range("G91").Value = "NA"
$ws.Range("H91").Value = "NA"
$ws.Range("I91").Value = "NA"
$ws.Range("J91").Value = "NA"
$ws.Range("K91").Value = "NA"
$ws.Range("L91").Value = $null
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = $null
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = $null
$ws.Range("D94").Value = -2280000
$ws.Range("E94").Value = -1414000
$ws.Range("F94").Value = -1780000
$ws.Range("G94").Value = "NA"
$ws.Range("H94").Value = "NA"
$ws.Range("I94").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("K94").Value = "NA"
$ws.Range("L94").Value = $null
$ws.Range("D95").Value = $null
$ws.Range("E95").Value = $null
$ws.Range("F95").Value = $null
$ws.Range("G95").Value = $null
$ws.Range("H95").Value = $null
$ws.Range("I95").Value = $null
$ws.Range("J95").Value = $null
$ws.Range("K95").Value = $null
$ws.Range("L95").Value = $null
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = -19000
$ws.Range("F96").Value = -19000
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = $null
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = $null
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = $null
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = $null
$ws.Range("D100").Value = 739000
$ws.Range("E100").Value = 1181000
$ws.Range("F100").Value = 1473000
$ws.Range("G100").Value = "NA"
$ws.Range("H100").Value = "NA"
$ws.Range("I100").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("K100").Value = "NA"
$ws.Range("L100").Value = $null
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = $null
$ws.Range("D102").Value = -644000
$ws.Range("E102").Value = 168000
$ws.Range("F102").Value = 130000
$ws.Range("G102").Value = "NA"
$ws.Range("H102").Value = "NA"
$ws.Range("I102").Value = "NA"
$ws.Range("J102").Value = "NA"
$ws.Range("K102").Value = "NA"
$ws.Range("L102").Value = $null
